$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Insert 7 new columns at H..N to hold the shared metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that every other property sheet already carries.
$ws.Columns("H:N").Insert()

# --- Header row (row 1) was wrongly populated with data; replace with the
# real field names, same as every other sheet (B=name, C=capacity, ...). ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2) ---
# C2 held the text "6cylinder"; split it into a numeric capacity value.
$ws.Range("C2").Value = 6
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-03-28"
$ws.Range("K2").Value = "陳唐山"
$ws.Range("L2").Value = 645
$ws.Range("M2").Value = "tmp38461"
$ws.Range("N2").Value = 33
